# Update Neural network Predictions (columns I and J) for rows 2-34
# with the redone scaling values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    @(2, 488909.5040597171, 61277.8806118276),
    @(3, 931710.2705830336, 263981.2768123746),
    @(4, 995722.7544250488, 889964.4019101858),
    @(5, 444743.693096593, -33794.666204134),
    @(6, 542559.0148231685, 85554.71745926328),
    @(7, 697175.7171912044, -231018.7562613264),
    @(8, 1121669.662679195, 553662.087174207),
    @(9, 1337418.48323679, 196037.4848187975),
    @(10, 576925.9032470435, -85463.24567181429),
    @(11, 842443.2461189777, 306978.2304192185),
    @(12, 4354507.453201294, 4618134.58324635),
    @(13, 4703065.762859344, 5555381.878884673),
    @(14, 930251.5227105319, 693680.2279697955),
    @(15, 868746.7246237099, 578109.2142707109),
    @(16, 824489.3397184461, -27162.14381982572),
    @(17, 3729374.387608766, 4173302.386969208),
    @(18, 692463.9237689823, -165414.2367077135),
    @(19, 573241.5970350653, 78910.80578934588),
    @(20, 503057.2439995557, 49654.73159999214),
    @(21, 427422.9901175051, -45016.4067395348),
    @(22, 1579029.358103961, 857438.5429856479),
    @(23, 932511.5042719841, 188271.0645669959),
    @(24, 1931934.901898682, 2410898.795298218),
    @(25, 1012476.85321334, 486977.1559833884),
    @(26, 1146408.416774929, 1045353.07312119),
    @(27, 773613.0433912724, 309692.1490680873),
    @(28, 405724.4347204789, 422300.2120257913),
    @(29, 471603.1526019573, -80155.08525077438),
    @(30, 3146365.948276877, 3296466.926583468),
    @(31, 1259196.869673312, 426429.7729101478),
    @(32, 927433.2107890844, 459193.2235685587),
    @(33, 5626236.92294538, 4720451.022147894),
    @(34, 2275496.022044897, 2505966.541249394)
)

foreach ($entry in $newValues) {
    $row = $entry[0]
    $ws.Cells.Item($row, 9).Value = $entry[1]
    $ws.Cells.Item($row, 10).Value = $entry[2]
}
